# Daily attendance processing - 2025-12-10 07:30:27
# For every "Recorded By" cell (column G) whose value starts with the
# literal prefix "System, " followed by additional recorder name(s),
# move "System" from the front of the comma-separated list to the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "System, "
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val.StartsWith($prefix)) {
        $rest = $val.Substring($prefix.Length)
        $newVal = $rest + ", System"
        $cell.Value = $newVal
    }
}
